# The original sheet has columns A:E (A,B,C,D,F headers) holding numeric
# measurements. The edit inserts a brand-new leading "ID" column that
# labels each data row, shifting the existing columns one place to the
# right (A->B, B->C, C->D, D->E, E->F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; this shifts all existing
# data/styles from columns A:E to B:F automatically.
$ws.Columns.Item(1).Insert()

# Give the new header cell (A1) the same look as the other header cells
# (bold font, thin border, centered/top alignment) by copying the format
# from the neighboring header cell B1, then set its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "ID"

# Row labels for the new ID column (rows 2-25).
$ids = @(
    "Hb 2",
    "Hb 3",
    "S 24",
    "S 28",
    "Hb 107",
    "Hb 66",
    "Hb 69",
    "Hb 95",
    "Hb 99",
    "Hb 92",
    "Hb 40",
    "Hb 41",
    "S 11",
    "Hb 57",
    "S 21",
    "S 22",
    "S 3",
    "S 4",
    "S 5",
    "Hb 74",
    "Hb 79",
    "Hb 32",
    "S 15",
    "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}

Write-Host "done"
